# Auto-generated edit script: apply cryptos price/coin updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'250.75"
$ws.Range("D3").Value = "'23.03"
$ws.Range("D4").Value = "'5.449"
$ws.Range("D5").Value = "'0.05673"
$ws.Range("D6").Value = "'3.417"
$ws.Range("D7").Value = "'6.387"
$ws.Range("D8").Value = "'0.8146"
$ws.Range("D9").Value = "'0.9201"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1439"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07509"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03110"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03089"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09376"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.553"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001593"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04764"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005791"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.006379"
$ws.Range("D20").Value = "'0.004996"
$ws.Range("D24").Value = "'2.204"
$ws.Range("D28").Value = "'0.0003001"
$ws.Range("D40").Value = "'0.04023"
$ws.Range("D41").Value = "'0.006768"
$ws.Range("D42").Value = "'0.1073"
$ws.Range("D44").Value = "'0.008035"
$ws.Range("D45").Value = "'0.00005804"
$ws.Range("D47").Value = "'0.5001"
$ws.Range("E48").Value = "47BOLOBOLOBestin24h"
$ws.Range("D49").Value = "'0.00002100"
